$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value  = 6.150599999999995
$ws.Range("B18").Value = 6.534799999999998
$ws.Range("B20").Value = 8.935200000000002
$ws.Range("B27").Value = 6.079400000000005
$ws.Range("B35").Value = 8.806500000000002
$ws.Range("B69").Value = 5.580399999999994
$ws.Range("B76").Value = 5.635499999999997
$ws.Range("B78").Value = 9.842699999999999
$ws.Range("B82").Value = 5.336600000000002
$ws.Range("B83").Value = 5.596999999999998
$ws.Range("B93").Value = 6.099299999999999
